$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update progress values (column C) ---
$ws.Range("C5").Value  = 0.9
$ws.Range("C6").Value  = 0.9
$ws.Range("C7").Value  = 1
$ws.Range("C8").Value  = 0.8
$ws.Range("C9").Value  = 0.4
$ws.Range("C10").Value = 0.6
$ws.Range("C11").Value = 0.5

$ws.Range("C13").Value = 1
$ws.Range("C14").Value = 0.5
$ws.Range("C15").Value = 0.3
$ws.Range("C16").Value = 0.5

$ws.Range("C18").Value = 0.5
$ws.Range("C19").Value = 0.8
$ws.Range("C20").Value = 0.6
$ws.Range("C21").Value = 0.3
$ws.Range("C22").Value = 0.3

$ws.Range("C24").Value = 0.5

# --- Add comments (column D) ---
# Order matches the order new strings were appended to the shared-string table.
$ws.Range("D5").Value  = "Regrouper les fichier en 1 fichier à analyser ?"
$ws.Range("D14").Value = "Revoir la standartisation avant la séparation pour l'implementer dans une pipline"
$ws.Range("D6").Value  = "Faire labelEncoder pour les variable binaire"
$ws.Range("D8").Value  = "Revoir la méthode peux être améliorer"

# --- Extend the global-average formula to include the newly tracked rows ---
$ws.Range("F6").Formula = "=AVERAGE(C5:C11,C13:C16,C18:C22,C24:C26)"

# --- Conditional formatting: drop the old C24:C28 color scale and rescope it to C24:C26 ---
$ws.Range("C24:C28").FormatConditions.Item(1).Delete()
$ws.Range("C24:C26").FormatConditions.AddColorScale(2) | Out-Null

# --- Sheet view: scroll so row 13 is at the top, and select C22 ---
$win = $excel.ActiveWindow
$win.ScrollRow = 13
$win.ScrollColumn = 1
$ws.Range("C22").Select() | Out-Null
